$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 10 so rows 10-17 shift down to 11-18,
# making room for the new scraped listing ahead of the existing ones.
$ws.Rows.Item(10).Insert()

$rows = @(
    @{ Row = 2; A = "2025-12-17 01:52:41"; B = "大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5455098"; G = 375; H = "🔥AI,Ai ◆開発" },
    @{ Row = 3; A = "2025-12-17 01:52:41"; B = "【フルリモート】官公庁向けPythonアプリ開発PM募集|7名チーム統括"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5454985"; G = 295; H = "🔥Python ◆開発 ◇アプリ" },
    @{ Row = 4; A = "2025-12-17 01:52:41"; B = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5217096"; G = 243; H = "🔥API ◆ツール" },
    @{ Row = 5; A = "2025-12-17 01:52:41"; B = "【Java/対話システム/心理学実験】協同問題解決プラットフォームの改修開発"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5439921"; G = 155; H = "★Java ◆開発" },
    @{ Row = 6; A = "2025-12-17 01:52:41"; B = "ホットペッパービューティーブログ一括投稿システム開発"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5455160"; G = 113; H = "◆開発,システム開発" },
    @{ Row = 7; A = "2025-12-17 01:52:41"; B = "【急募】新規システム開発に伴う要件定義依頼"; C = "システム開発"; D = "10,000 円 ~ 20,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5455415"; G = 110; H = "◆開発,システム開発" },
    @{ Row = 8; A = "2025-12-17 01:52:41"; B = "【急募】iPhone・Android対応の天気アプリ開発をお願いします!"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5455038"; G = 100; H = "◆開発 ◇アプリ" },
    @{ Row = 9; A = "2025-12-17 01:52:41"; B = "Kabuステーション自動売買アプリの開発"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5455251"; G = 93; H = "◆開発 ◇アプリ" },
    @{ Row = 10; A = "2025-12-17 01:52:41"; B = "【急募】MVNO会員向けマイページ新規開発エンジニア募集"; C = "システム開発"; D = "1,000,000 円 ~ 3,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5455513"; G = 75; H = "◆開発" },
    @{ Row = 11; A = "2025-12-17 01:52:41"; B = "ホームページ診断チェックツール"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5455029"; G = 73; H = "◆ツール" },
    @{ Row = 12; A = "2025-12-17 01:52:41"; B = "【急募】帳票デジタル化のフロントエンド開発者募集"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5454857"; G = 68; H = "◆開発" },
    @{ Row = 13; A = "2025-12-17 01:52:41"; B = "【急募】Accessシステム改修・CSV読込・MySQLクラウド化・PDFデータ調整"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5455015"; G = 53; H = "◇MySQL" },
    @{ Row = 14; A = "2025-12-17 01:52:41"; B = "【改善提案募集】事業管理スプレッドシートの見直し・改善提案をお願いします。"; C = "システム開発"; D = "1,000 ~ 5,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5455422"; G = 30; H = "◇管理" },
    @{ Row = 15; A = "2025-12-17 01:52:41"; B = "【急募】wixシステムでのメッセージ送信システム構築依頼"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5455067"; G = 33; H = $null },
    @{ Row = 16; A = "2025-12-17 01:52:41"; B = "【急募】企業のセキュリティ対策を担うエンジニア募集"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5450345"; G = 25; H = $null },
    @{ Row = 17; A = "2025-12-17 01:52:41"; B = "〖リモート可〗Delphiエンジニア募集"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5341051"; G = 25; H = $null },
    @{ Row = 18; A = "2025-12-17 01:52:41"; B = "【SESエンジニア募集】多様なプロジェクトに参画可能!"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5437544"; G = 25; H = $null }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    if ($item.H -ne $null) {
        $ws.Cells.Item($r, 8).Value = $item.H
    } else {
        $ws.Cells.Item($r, 8).ClearContents()
    }
}

# Hyperlinks do not automatically re-target after Rows.Insert(), so rebuild them
# from scratch in row order (F2..F18) to match the shifted URLs.
$ws.Hyperlinks.Delete()
foreach ($item in $rows) {
    $ws.Hyperlinks.Add($ws.Cells.Item($item.Row, 6), $item.F)
}

# Column D widened from 30 to 32 characters. ColumnWidth (character units) needs
# an input a bit under 32 because Excel rounds up to the nearest pixel internally;
# 31.14 is the value that lands exactly on width="32" in the saved file.
$ws.Columns.Item(4).ColumnWidth = 31.14

Write-Host "edit complete"
